# "avances en carga masiva"
# Adds the REP_* / CNAT_* (representative / natural-person client) header
# columns to the "DATA" bulk-upload template sheet, extending the table
# from columns A:Y out to A:AP, and updates the sheet view (zoom + the
# remembered selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, in column order starting right after NUMERO_TELEFONO_2 (Y1).
$headers = @(
  "REP_DOC_ID",
  "REP_PRIMER_NOMBRE",
  "REP_SEGUNDO_NOMBRE",
  "REP_PRIMER_APELLIDO",
  "REP_SEGUNDO_APELLIDO",
  "REP_EMAIL",
  "REP_RIF",
  "CNAT_DOC_ID",
  "CNAT_PRIMER_NOMBRE",
  "CNAT_SEGUNDO_NOMBRE",
  "CNAT_PRIMER_APELLIDO",
  "CNAT_SEGUNDO_APELLIDO",
  "CNAT_FECHA_NACIMIENTO",
  "CNAT_TIPO_NATURAL",
  "CNAT_GENERO",
  "CNAT_PROFESION",
  "CNAT_T_DOC"
)

$startCol = 26   # column Z

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $startCol + $i).Value = $headers[$i]
}

# Column widths (bestFit-style) for the new columns, matching the
# widths the workbook carries after adding the header text. Column Y
# (25) picks up the same bestFit width as X (24); the rest get their own.
$widths = @{
  25 = 20.592447916666668
  26 = 11.166666666666666
  27 = 20.166666666666668
  28 = 22.307291666666668
  29 = 20.877604166666668
  30 = 22.877604166666668
  31 = 9.877604166666666
  32 = 7.166666666666667
  33 = 12.736979166666666
  34 = 21.877604166666668
  35 = 23.877604166666668
  36 = 22.451822916666668
  37 = 24.451822916666668
  38 = 24.592447916666668
  39 = 19.592447916666668
  40 = 13.307291666666666
  41 = 16.307291666666668
  42 = 11.877604166666666
}

foreach ($col in $widths.Keys) {
  $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}

# Sheet view: zoom to 110% and move the remembered selection to AQ3
# (just past the new last used column, AP).
$ws.Application.ActiveWindow.Zoom = 110
$ws.Range("AQ3").Select() | Out-Null
